$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.884.27'
$ws.Range("E2").Value = '  +4.83%  '
$ws.Range("D3").Value = '2.265.44'
$ws.Range("E3").Value = '  +2.05%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '301.95'
$ws.Range("E5").Value = '  +3.35%  '
$ws.Range("D6").Value = '92.25'
$ws.Range("E6").Value = '  +6.16%  '
$ws.Range("E7").Value = '  +3.38%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '0.485'
$ws.Range("E9").Value = '  +3.91%  '
$ws.Range("E10").Value = '  +8.19%  '
$ws.Range("D11").Value = '32.28'
$ws.Range("E11").Value = '  +6.03%  '
$ws.Range("D12").Value = '0.0798'
$ws.Range("E12").Value = '  +2.32%  '
$ws.Range("E13").Value = '  +2.08%  '
$ws.Range("E14").Value = '  +3.58%  '
$ws.Range("D15").Value = '2.615.68'
$ws.Range("E15").Value = '  +2.09%  '
$ws.Range("D16").Value = '14.17'
$ws.Range("E16").Value = '  +2.80%  '
$ws.Range("D17").Value = '2.270.82'
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("E18").Value = '  +3.58%  '
$ws.Range("D19").Value = '41.777.40'
$ws.Range("E19").Value = '  +4.78%  '
$ws.Range("E20").Value = '  +8.39%  '
$ws.Range("D21").Value = '0.0₃0904'
$ws.Range("E21").Value = '  +2.00%  '
$ws.Range("E22").Value = '  +3.47%  '
$ws.Range("D23").Value = '66.96'
$ws.Range("D24").Value = '241.84'
$ws.Range("E24").Value = '  +2.09%  '
$ws.Range("E25").Value = '  +3.90%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  +4.01%  '
$ws.Range("D28").Value = '23.87'
$ws.Range("E28").Value = '  +2.64%  '
$ws.Range("D29").Value = '9.63'
$ws.Range("E29").Value = '  +4.41%  '
$ws.Range("E30").Value = '  -12.05%  '
$ws.Range("D31").Value = '159.28'
$ws.Range("E31").Value = '  +1.06%  '
$ws.Range("D32").Value = '33.77'
$ws.Range("E32").Value = '  +6.20%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").Value = '5.17'
$ws.Range("E34").Value = '  +3.97%  '
$ws.Range("E35").Value = '  +4.47%  '
$ws.Range("E36").Value = '  +3.14%  '
$ws.Range("E37").Value = '  +2.07%  '
$ws.Range("E38").Value = '  +5.84%  '
$ws.Range("E39").Value = '  +3.49%  '
$ws.Range("E40").Value = '  +8.87%  '
$ws.Range("E41").Value = '  +4.63%  '
$ws.Range("E42").Value = '  +5.57%  '
$ws.Range("D43").Value = '2.074.86'
$ws.Range("E43").Value = '  -0.65%  '
$ws.Range("D44").Value = '19.57'
$ws.Range("E44").Value = '  +9.33%  '
$ws.Range("E45").Value = '  +3.28%  '
$ws.Range("D46").Value = '10.17'
$ws.Range("E46").Value = '  +3.78%  '
$ws.Range("E47").Value = '  +8.54%  '
$ws.Range("E48").Value = '  +2.37%  '
$ws.Range("E49").Value = '  +3.62%  '
$ws.Range("E50").Value = '  +3.20%  '
$ws.Range("D51").Value = '51.77'
$ws.Range("E51").Value = '  +5.73%  '
